# Generate Report for Handoff
#
# cdc20666-6a44-42da-83ab-4a0b3709cda3 just received a new handoff event,
# so its row moves to the bottom of its block (rows 6-8) while
# e44d024b-48d5-43b9-8fcd-23440750096a and 878c6d98-d1ce-4a46-b911-cee247c17219
# shift up one row, on every sheet (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: columns A (File Name), B (zh-cn status), C (de-de status),
# D (Latest Handoff Date)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "e44d024b-48d5-43b9-8fcd-23440750096a.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "In Translation"
$ws.Range("D6").Value = "2016-30-12 14:30:16"

$ws.Range("A7").Value = "878c6d98-d1ce-4a46-b911-cee247c17219.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "2016-29-12 14:29:09"

$ws.Range("A8").Value = "cdc20666-6a44-42da-83ab-4a0b3709cda3.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-33-12 14:33:35"

# ---------------------------------------------------------------------
# zh-cn sheet: columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A6").Value = "e44d024b-48d5-43b9-8fcd-23440750096a.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("D6").Value = "e44d024b-48d5-43b9-8fcd-23440750096a.0336dd5674b37fb5d2d98fd7fb595158049c241b.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-12 14:30:09"

$ws.Range("A7").Value = "878c6d98-d1ce-4a46-b911-cee247c17219.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "878c6d98-d1ce-4a46-b911-cee247c17219.7112a32ca86eda0539d4e496fef1abfc21660fa3.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-12 14:29:06"

$ws.Range("A8").Value = "cdc20666-6a44-42da-83ab-4a0b3709cda3.md"
$ws.Range("B8").Value = ".md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "cdc20666-6a44-42da-83ab-4a0b3709cda3.8b4f9242a0c1edb75fa34c4d673ff4550f445d49.zh-cn.xlf"
$ws.Range("E8").Value = "2016-03-12 14:33:32"

$ws.Range("A9").Value = "ff71d0f1-b8a4-4430-8039-78eaf535d27b.md"
$ws.Range("B9").Value = ".md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "ff71d0f1-b8a4-4430-8039-78eaf535d27b.89a82442c4c7ee1a54185e62364695a3c2e7c2d6.zh-cn.xlf"
$ws.Range("E9").Value = "2016-03-12 14:31:03"

# ---------------------------------------------------------------------
# de-de sheet: columns A (Source File Name), B (File Extension),
# C (Status), D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A6").Value = "e44d024b-48d5-43b9-8fcd-23440750096a.md"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "In Translation"
$ws.Range("D6").Value = "e44d024b-48d5-43b9-8fcd-23440750096a.0336dd5674b37fb5d2d98fd7fb595158049c241b.de-de.xlf"
$ws.Range("E6").Value = "2016-03-12 14:30:16"

$ws.Range("A7").Value = "878c6d98-d1ce-4a46-b911-cee247c17219.md"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "878c6d98-d1ce-4a46-b911-cee247c17219.7112a32ca86eda0539d4e496fef1abfc21660fa3.de-de.xlf"
$ws.Range("E7").Value = "2016-03-12 14:29:09"

$ws.Range("A8").Value = "cdc20666-6a44-42da-83ab-4a0b3709cda3.md"
$ws.Range("B8").Value = ".md"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "cdc20666-6a44-42da-83ab-4a0b3709cda3.8b4f9242a0c1edb75fa34c4d673ff4550f445d49.de-de.xlf"
$ws.Range("E8").Value = "2016-03-12 14:33:35"

$ws.Range("A9").Value = "ff71d0f1-b8a4-4430-8039-78eaf535d27b.md"
$ws.Range("B9").Value = ".md"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "ff71d0f1-b8a4-4430-8039-78eaf535d27b.89a82442c4c7ee1a54185e62364695a3c2e7c2d6.de-de.xlf"
$ws.Range("E9").Value = "2016-03-12 14:31:07"
